# Generate Report for Handoff
# Updates the localization-status report: marks rows 7,8,9,10,13,14 (the
# "Ready for handoff" entries) as handoff-type "ht" on the zh-cn and de-de
# sheets, and refreshes the "Latest Handoff"/"Latest HO Xliff Generate Date"
# timestamps to reflect the new handoff generation run.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 10, 13, 14)

# --- Overview sheet: refresh "Latest HO Xliff Generate Date" (column G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-12 12:22:24"
}

# --- zh-cn sheet: set Priority (column E) to "ht" and refresh
#     "Latest Handoff Datetime" (column H) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-12 12:22:17"
}

# --- de-de sheet: set Priority (column E) to "ht" and refresh
#     "Latest Handoff Datetime" (column H) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-12 12:22:24"
}
